$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin name / Link columns (row reordering for TRON/WrappedEther, Filecoin/RenderToken, VeChain/Hedera)
$ws.Range("B12").Value2 = 'WrappedEther'
$ws.Range("C12").Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("B13").Value2 = 'TRON'
$ws.Range("C13").Value2 = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("B34").Value2 = 'Filecoin'
$ws.Range("C34").Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("B35").Value2 = 'RenderToken'
$ws.Range("C35").Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("B37").Value2 = 'Hedera'
$ws.Range("C37").Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("B38").Value2 = 'VeChain'
$ws.Range("C38").Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

# Update Price column (column D) -- force text so numeric-looking strings are preserved verbatim
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = '26.846.96'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = '1.838.34'
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = '1.007'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '308.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '1.005'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '0.4719'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '0.3646'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '0.07143'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '0.9162'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '19.47'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '1.899.20'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '0.07595'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '5.273'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '6.382'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '87.71'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = '1.008'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '0.000008620'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '1.005'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '26.872.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = '14.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '5.006'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '10.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '1.928'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '151.43'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '18.16'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '2.000'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '113.98'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '4.851'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '0.08815'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = '3.226'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '1.165'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '0.7392'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '4.467'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = '2.741'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '1.086'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '0.05241'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '0.01938'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '2.966'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '0.5163'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '6.922'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '0.1508'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '8.128'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '10.40'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '0.4690'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '1.006'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '101.41'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '1.589'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '65.25'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '0.06032'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = '0.8824'
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) column (column E)
$ws.Range("E2").Value2 = '  +1.11%  '
$ws.Range("E3").Value2 = '  +1.48%  '
$ws.Range("E4").Value2 = '  +0.38%  '
$ws.Range("E5").Value2 = '  +0.96%  '
$ws.Range("E6").Value2 = '  +0.23%  '
$ws.Range("E7").Value2 = '  +3.73%  '
$ws.Range("E8").Value2 = '  +1.42%  '
$ws.Range("E9").Value2 = '  +0.46%  '
$ws.Range("E10").Value2 = '  +2.60%  '
$ws.Range("E11").Value2 = '  +0.81%  '
$ws.Range("E12").Value2 = '  +6.58%  '
$ws.Range("E13").Value2 = '  -1.57%  '
$ws.Range("E14").Value2 = '  +0.33%  '
$ws.Range("E15").Value2 = '  +1.38%  '
$ws.Range("E16").Value2 = '  +1.02%  '
$ws.Range("E17").Value2 = '  +0.27%  '
$ws.Range("E18").Value2 = '  +0.80%  '
$ws.Range("E19").Value2 = '  +0.24%  '
$ws.Range("E20").Value2 = '  +1.10%  '
$ws.Range("E21").Value2 = '  +2.42%  '
$ws.Range("E22").Value2 = '  +0.83%  '
$ws.Range("E23").Value2 = '  +0.64%  '
$ws.Range("E24").Value2 = '  +0.00%  '
$ws.Range("E25").Value2 = '  -0.17%  '
$ws.Range("E26").Value2 = '  +2.04%  '
$ws.Range("E27").Value2 = '  -0.83%  '
$ws.Range("E28").Value2 = '  +1.34%  '
$ws.Range("E29").Value2 = '  +0.35%  '
$ws.Range("E30").Value2 = '  +0.97%  '
$ws.Range("E31").Value2 = '  +3.06%  '
$ws.Range("E32").Value2 = '  +4.88%  '
$ws.Range("E33").Value2 = '  +0.09%  '
$ws.Range("E34").Value2 = '  +0.79%  '
$ws.Range("E35").Value2 = '  +0.97%  '
$ws.Range("E36").Value2 = '  +1.43%  '
$ws.Range("E37").Value2 = '  +3.28%  '
$ws.Range("E38").Value2 = '  +0.12%  '
$ws.Range("E39").Value2 = '  +1.78%  '
$ws.Range("E40").Value2 = '  +1.52%  '
$ws.Range("E41").Value2 = '  +2.16%  '
$ws.Range("E42").Value2 = '  +0.14%  '
$ws.Range("E43").Value2 = '  +1.61%  '
$ws.Range("E44").Value2 = '  +4.13%  '
$ws.Range("E45").Value2 = '  +0.02%  '
$ws.Range("E46").Value2 = '  +0.24%  '
$ws.Range("E47").Value2 = '  +2.18%  '
$ws.Range("E48").Value2 = '  +1.31%  '
$ws.Range("E49").Value2 = '  +2.53%  '
$ws.Range("E50").Value2 = '  +0.51%  '
$ws.Range("E51").Value2 = '  +4.30%  '

Write-Host "Applied cryptos list update"
